$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Home")

# --- Paste the debug console log (home.js) into column G, next to the existing
# Home_N_*/Home_S_* scene-name table, mirroring the values used while tracking down
# the if/else bug in fetchMode(). Label rows get a light-green highlight fill.

$ws.Range("G3").Value = 'Home_N_Main_Fb'
$ws.Range("G3").Interior.Color = 14348258
$ws.Range("G4").Value = 'home.js:497 undefined'
$ws.Range("G5").Value = 'home.js:496 Home_S_main_Fb'
$ws.Range("G5").Interior.Color = 14348258
$ws.Range("G6").Value = 'home.js:497 undefined'
$ws.Range("G7").Value = 'home.js:496 Home_N_W_St_Fb'
$ws.Range("G7").Interior.Color = 14348258
$ws.Range("G8").Value = 'home.js:497 undefined'
$ws.Range("G9").Value = 'home.js:496 Home_S_W_St_Fb'
$ws.Range("G9").Interior.Color = 14348258
$ws.Range("G10").Value = 'home.js:497 undefined'
$ws.Range("G11").Value = 'home.js:496 Home_N_W_Puls_Fb'
$ws.Range("G11").Interior.Color = 14348258
$ws.Range("G12").Value = 'home.js:497 undefined'
$ws.Range("G13").Value = 'home.js:496 Home_S_W_Puls_Fb'
$ws.Range("G13").Interior.Color = 14348258
$ws.Range("G14").Value = 'home.js:497 undefined'
$ws.Range("G15").Value = 'home.js:496 Home_N_W_St+Puls_Fb'
$ws.Range("G15").Interior.Color = 14348258
$ws.Range("G16").Value = 'home.js:497 undefined'
$ws.Range("G17").Value = 'home.js:496 Home_S_W_St+Puls_Fb'
$ws.Range("G17").Interior.Color = 14348258
$ws.Range("G18").Value = 'home.js:497 127'
$ws.Range("G19").Value = 'home.js:496 Home_N_sunset_Fb'
$ws.Range("G19").Interior.Color = 14348258
$ws.Range("G20").Value = 'home.js:497 undefined'
$ws.Range("G21").Value = 'home.js:496 Home_S_Color_dynamic_Fb'
$ws.Range("G21").Interior.Color = 14348258
$ws.Range("G22").Value = 'home.js:497 undefined'
$ws.Range("G23").Value = 'home.js:496 Home_N_autumn_Fb'
$ws.Range("G23").Interior.Color = 14348258
$ws.Range("G24").Value = 'home.js:497 undefined'
$ws.Range("G25").Value = 'home.js:496 Home_S_ECO_Fb'
$ws.Range("G25").Interior.Color = 14348258
$ws.Range("G26").Value = 'home.js:497 undefined'
$ws.Range("G27").Value = 'home.js:496 Home_N_ECO_Fb'
$ws.Range("G27").Interior.Color = 14348258
$ws.Range("G28").Value = 'home.js:497 undefined'
$ws.Range("G29").Value = 'home.js:496 Home_S_Violet_Fb'
$ws.Range("G29").Interior.Color = 14348258
$ws.Range("G30").Value = 'home.js:497 undefined'
$ws.Range("G31").Value = 'home.js:496 Home_N_Violet_Fb'
$ws.Range("G31").Interior.Color = 14348258
$ws.Range("G32").Value = 'home.js:497 undefined'
$ws.Range("G33").Value = 'home.js:496 Home_S_lightbox_Fb'
$ws.Range("G34").Value = 'home.js:497 144'
$ws.Range("G35").Value = 'home.js:496 Home_N_lightbox_Fb'
$ws.Range("G36").Value = 'home.js:497 undefined'
$ws.Range("G37").Value = 'home.js:496 Home_S_logotypes_Fb'
$ws.Range("G38").Value = 'home.js:497 145'
$ws.Range("G39").Value = 'home.js:496 Home_N_logotypes_Fb'
$ws.Range("G40").Value = 'home.js:497 undefined'

# New column is wide enough to show the full console-log text
$ws.Columns.Item(7).ColumnWidth = 40.109375

# Leave the selection where it ended up after pasting the log (matches the saved file)
$ws.Range("H41").Select()
